$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of cell -> new text value (mirrors the source diff).
$updates = [ordered]@{
    'D2' = '42.493.87'
    'E2' = '  -1.51%  '
    'D3' = '2.540.38'
    'E3' = '  -0.55%  '
    'E4' = '  -0.08%  '
    'D5' = '312.11'
    'E5' = '  -1.57%  '
    'D6' = '99.04'
    'E6' = '  +1.86%  '
    'D7' = '0.568'
    'E7' = '  -1.13%  '
    'E8' = '  +0.03%  '
    'E9' = '  -2.98%  '
    'D10' = '35.62'
    'E10' = '  -0.29%  '
    'E11' = '  -1.26%  '
    'D12' = '7.34'
    'E12' = '  -2.05%  '
    'E13' = '  +0.05%  '
    'D14' = '2.931.65'
    'E14' = '  -0.53%  '
    'D15' = '16.13'
    'E15' = '  +7.22%  '
    'D16' = '2.591.74'
    'E16' = '  +0.67%  '
    'D17' = '0.837'
    'E17' = '  -1.44%  '
    'D18' = '42.487.95'
    'E18' = '  -1.62%  '
    'D19' = '6.76'
    'E19' = '  -1.06%  '
    'B20' = 'ShibaInu'
    'C20' = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
    'D20' = '0.0₃0947'
    'E20' = '  -1.78%  '
    'B21' = 'InternetComputer(DFINITY)'
    'C21' = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
    'D21' = '12.24'
    'E21' = '  -2.87%  '
    'D22' = '68.68'
    'E22' = '  -2.13%  '
    'D23' = '242.38'
    'E23' = '  -4.32%  '
    'D24' = '2.90'
    'E24' = '  -1.99%  '
    'E25' = '  -0.50%  '
    'E26' = '  +0.03%  '
    'D27' = '26.23'
    'E27' = '  -1.97%  '
    'E28' = '  -3.46%  '
    'D29' = '39.76'
    'E29' = '  -2.75%  '
    'D30' = '10.08'
    'E30' = '  -1.87%  '
    'D31' = '158.75'
    'E31' = '  +1.82%  '
    'D32' = '5.67'
    'E32' = '  -2.95%  '
    'D33' = '2.80'
    'E33' = '  +14.77%  '
    'D34' = '0.0796'
    'E34' = '  -0.60%  '
    'E35' = '  -3.16%  '
    'D36' = '2.04'
    'E36' = '  -3.97%  '
    'E37' = '  -4.94%  '
    'D38' = '17.99'
    'E38' = '  -6.88%  '
    'E39' = '  -1.42%  '
    'D40' = '0.118'
    'E40' = '  -0.57%  '
    'D41' = '4.18'
    'E41' = '  +7.91%  '
    'D42' = '21.73'
    'E42' = '  -1.37%  '
    'E43' = '  +0.10%  '
    'D44' = '3.31'
    'E44' = '  +2.03%  '
    'E45' = '  -2.94%  '
    'D46' = '1.953.84'
    'E46' = '  -1.73%  '
    'E47' = '  -1.66%  '
    'D48' = '2.786.82'
    'E48' = '  -0.64%  '
    'D49' = '80.51'
    'E49' = '  -4.94%  '
    'D50' = '0.192'
    'E50' = '  -0.53%  '
    'D51' = '72.34'
    'E51' = '  -2.83%  '
}

foreach ($cell in $updates.Keys) {
    $range = $ws.Range($cell)
    # Column D holds price text such as "42.493.87" or "0.568" that Excel
    # would otherwise auto-coerce into a number (dropping separators/
    # trailing zeros). Force text format so the literal string is kept,
    # matching the workbook's stored inline-string cells.
    if ($cell -like "D*") {
        $range.NumberFormat = "@"
    }
    $range.Value = $updates[$cell]
}